# Generate Report for handback
# - Updates status text "Ready for handoff" -> "Handed back: in sync with en-US"
#   (shared across Overview/zh-cn/de-de sheets)
# - Fills in "Latest Target File" / "Latest Handback File" columns (E/F) for the
#   two real source rows on the "zh-cn" and "de-de" sheets, with hyperlinks that
#   mirror the existing handoff-file hyperlinks
# - Stamps "Latest Handback DateTime" (G) with the handback timestamp

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: refresh the status text wherever it currently reads
# "Ready for handoff" (rows for the two real files, zh-cn + de-de columns).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

# Row 2 -> 222fddb2-7a44-47a2-998c-bd4d2d6e7aa4
$wsZh.Range("E2").Value = "222fddb2-7a44-47a2-998c-bd4d2d6e7aa4.md"
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/24abf57586c5cdbbb8d48461f430d1353fdc3ca3/e2e/222fddb2-7a44-47a2-998c-bd4d2d6e7aa4.md", "", "", "222fddb2-7a44-47a2-998c-bd4d2d6e7aa4.md")

$wsZh.Range("F2").Value = "222fddb2-7a44-47a2-998c-bd4d2d6e7aa4.dd16e3dfc077f2cb539d74d68aa419308c277e1d.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e6344398c3ea7eed37e77095ca104b8a8864bf0d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/222fddb2-7a44-47a2-998c-bd4d2d6e7aa4.dd16e3dfc077f2cb539d74d68aa419308c277e1d.zh-cn.xlf", "", "", "222fddb2-7a44-47a2-998c-bd4d2d6e7aa4.dd16e3dfc077f2cb539d74d68aa419308c277e1d.zh-cn.xlf")

$wsZh.Range("G2").Value = "2016-01-22 02:29:30"

# Row 3 -> 56bec4ee-7eef-4f4f-8950-234bdbffa32a
$wsZh.Range("E3").Value = "56bec4ee-7eef-4f4f-8950-234bdbffa32a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/24abf57586c5cdbbb8d48461f430d1353fdc3ca3/e2e/56bec4ee-7eef-4f4f-8950-234bdbffa32a.md", "", "", "56bec4ee-7eef-4f4f-8950-234bdbffa32a.md")

$wsZh.Range("F3").Value = "56bec4ee-7eef-4f4f-8950-234bdbffa32a.91bdf8ed8f7d8d162841005f9b9262c41dad8bec.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e6344398c3ea7eed37e77095ca104b8a8864bf0d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/56bec4ee-7eef-4f4f-8950-234bdbffa32a.91bdf8ed8f7d8d162841005f9b9262c41dad8bec.zh-cn.xlf", "", "", "56bec4ee-7eef-4f4f-8950-234bdbffa32a.91bdf8ed8f7d8d162841005f9b9262c41dad8bec.zh-cn.xlf")

$wsZh.Range("G3").Value = "2016-01-22 02:29:30"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

# Row 2 -> 222fddb2-7a44-47a2-998c-bd4d2d6e7aa4
$wsDe.Range("E2").Value = "222fddb2-7a44-47a2-998c-bd4d2d6e7aa4.md"
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/24abf57586c5cdbbb8d48461f430d1353fdc3ca3/e2e/222fddb2-7a44-47a2-998c-bd4d2d6e7aa4.md", "", "", "222fddb2-7a44-47a2-998c-bd4d2d6e7aa4.md")

$wsDe.Range("F2").Value = "222fddb2-7a44-47a2-998c-bd4d2d6e7aa4.dd16e3dfc077f2cb539d74d68aa419308c277e1d.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3295e1f39c4711270478a1054fb8555019b6db10/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/222fddb2-7a44-47a2-998c-bd4d2d6e7aa4.dd16e3dfc077f2cb539d74d68aa419308c277e1d.de-de.xlf", "", "", "222fddb2-7a44-47a2-998c-bd4d2d6e7aa4.dd16e3dfc077f2cb539d74d68aa419308c277e1d.de-de.xlf")

$wsDe.Range("G2").Value = "2016-01-22 02:29:53"

# Row 3 -> 56bec4ee-7eef-4f4f-8950-234bdbffa32a
$wsDe.Range("E3").Value = "56bec4ee-7eef-4f4f-8950-234bdbffa32a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/24abf57586c5cdbbb8d48461f430d1353fdc3ca3/e2e/56bec4ee-7eef-4f4f-8950-234bdbffa32a.md", "", "", "56bec4ee-7eef-4f4f-8950-234bdbffa32a.md")

$wsDe.Range("F3").Value = "56bec4ee-7eef-4f4f-8950-234bdbffa32a.91bdf8ed8f7d8d162841005f9b9262c41dad8bec.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3295e1f39c4711270478a1054fb8555019b6db10/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/56bec4ee-7eef-4f4f-8950-234bdbffa32a.91bdf8ed8f7d8d162841005f9b9262c41dad8bec.de-de.xlf", "", "", "56bec4ee-7eef-4f4f-8950-234bdbffa32a.91bdf8ed8f7d8d162841005f9b9262c41dad8bec.de-de.xlf")

$wsDe.Range("G3").Value = "2016-01-22 02:29:53"

Write-Host "Handback report generated."
